$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# workbook-level cosmetic metadata (author path + window geometry)
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Top = -110
$wb.Windows.Item(1).Left = -110
$wb.Windows.Item(1).Width = 19420
$wb.Windows.Item(1).Height = 10420

# ---------------------------------------------------------------------------
# GUI sheet ("testsheet2" / tab name "GUI") — add the two new columns that
# capture "Adults Count" and "Child Count" between "Check Out Date" and
# "Expected Hotel Name".
# ---------------------------------------------------------------------------
$gui = $wb.Worksheets.Item("GUI")

# Insert two new blank columns at E so the existing "Expected Hotel Name"
# column (and its data) shifts from E -> G, matching Excel's own behaviour
# when a user selects E:F and inserts columns.
$gui.Range("E1:F1").EntireColumn.Insert()

# New header cells.
$gui.Cells.Item(1, 5).Value = "Adults Count"
$gui.Cells.Item(1, 6).Value = "Child Count"

# New data values (plain numbers, not text).
$gui.Cells.Item(2, 5).Value = 4
$gui.Cells.Item(2, 6).Value = 2

# Match the left-aligned "normal" style used by the other plain-value cells
# in row 2 (style index 1), instead of the inherited header/border style.
$gui.Cells.Item(2, 5).Style = $gui.Cells.Item(2, 1).Style
$gui.Cells.Item(2, 6).Style = $gui.Cells.Item(2, 1).Style

# Selection / zoom, mirroring the saved view state after editing.
$gui.Application.ActiveWindow.Zoom = 130
$gui.Range("F1").Select()

# ---------------------------------------------------------------------------
# API sheet — only the view/zoom/selection changed (no data edits).
# ---------------------------------------------------------------------------
$api = $wb.Worksheets.Item("API")
$api.Select()
$api.Application.ActiveWindow.Zoom = 115
$api.Range("G1").Select()
